$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 75000
$ws.Range("J3").Value = 75000
$ws.Range("L3").Value = 75000
$ws.Range("N3").Value = -75228
$ws.Range("H17").Value = 347677.62
$ws.Range("J17").Value = 347677.62
$ws.Range("L17").Value = 1043032.86
$ws.Range("N17").Value = -1043368.86
$ws.Range("H102").Value = 75000
$ws.Range("J102").Value = 75000
$ws.Range("L102").Value = 75000
$ws.Range("N102").Value = -81490
$ws.Range("H121").Value = 4498.3335
$ws.Range("J121").Value = 6500
$ws.Range("L121").Value = 19500
$ws.Range("N121").Value = -22994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2102.6667
$ws.Range("I74").Value = 1536.5151
$ws.Range("J74").Value = 4178.5557
$ws.Range("K74").Value = 1536.5151
$ws.Range("L74").Value = 4178.5557
$ws.Range("M74").Value = -662.5151000000001
$ws.Range("N74").Value = -5926.5557
$ws.Range("H77").Value = 2102.6667
$ws.Range("I77").Value = 1536.5151
$ws.Range("J77").Value = 4178.5557
$ws.Range("K77").Value = 7682.575500000001
$ws.Range("L77").Value = 20892.7785
$ws.Range("M77").Value = -3314.575500000001
$ws.Range("N77").Value = -29628.7785
$ws.Range("H93").Value = 68741.336
$ws.Range("J93").Value = 68741.336
$ws.Range("L93").Value = 68741.336
$ws.Range("N93").Value = -73733.336
$ws.Range("H101").Value = 72249.5
$ws.Range("J101").Value = 72249.5
$ws.Range("L101").Value = 72249.5
$ws.Range("N101").Value = -78739.5
$ws.Range("H122").Value = 4535.3335
$ws.Range("I122").Value = 4487.4546
$ws.Range("K122").Value = 13462.3638
$ws.Range("M122").Value = -11012.3638
$ws.Range("H132").Value = 4008.0334
$ws.Range("I132").Value = 1416.4073
$ws.Range("J132").Value = 27332.666
$ws.Range("K132").Value = 4249.2219
$ws.Range("L132").Value = 81997.99800000001
$ws.Range("M132").Value = -1719.2219
$ws.Range("N132").Value = -87057.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3041.3333
$ws.Range("I94").Value = 2393.1333
$ws.Range("J94").Value = 3851.5833
$ws.Range("K94").Value = 2393.1333
$ws.Range("L94").Value = 3851.5833
$ws.Range("M94").Value = -1942.1333
$ws.Range("N94").Value = -4753.5833
$ws.Range("H99").Value = 3582.7334
$ws.Range("I99").Value = 2254.2
$ws.Range("K99").Value = 2254.2
$ws.Range("M99").Value = -756.1999999999998
$ws.Range("H105").Value = 2697.7942
$ws.Range("I105").Value = 1998.2413
$ws.Range("K105").Value = 1998.2413
$ws.Range("M105").Value = -251.2412999999999
$ws.Range("H134").Value = 4502.6104
$ws.Range("I134").Value = 1298.2
$ws.Range("J134").Value = 9008.8125
$ws.Range("K134").Value = 3894.6
$ws.Range("L134").Value = 27026.4375
$ws.Range("M134").Value = -1359.6
$ws.Range("N134").Value = -32096.4375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H51").Value = 50783.332
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9264
$ws.Range("H59").Value = 93850
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = ""
$ws.Range("H60").Value = 9800
$ws.Range("I60").Value = 9800
$ws.Range("K60").Value = 9800
$ws.Range("M60").Value = -9289
$ws.Range("H61").Value = 50783.332
$ws.Range("I61").Value = 10000
$ws.Range("K61").Value = 10000
$ws.Range("M61").Value = -9652
$ws.Range("H86").Value = 7089.6
$ws.Range("I86").Value = 6737.125
$ws.Range("J86").Value = 8499.5
$ws.Range("K86").Value = 6737.125
$ws.Range("L86").Value = 8499.5
$ws.Range("M86").Value = -5614.125
$ws.Range("N86").Value = -10745.5
$ws.Range("H89").Value = 7089.6
$ws.Range("I89").Value = 6737.125
$ws.Range("J89").Value = 8499.5
$ws.Range("K89").Value = 33685.625
$ws.Range("L89").Value = 42497.5
$ws.Range("M89").Value = -28069.625
$ws.Range("N89").Value = -53729.5
$ws.Range("H107").Value = 2688.0386
$ws.Range("I107").Value = 1980.7222
$ws.Range("K107").Value = 1980.7222
$ws.Range("M107").Value = -60.72219999999993
$ws.Range("H122").Value = 3486.861
$ws.Range("I122").Value = 3532.8462
$ws.Range("J122").Value = 3367.3
$ws.Range("K122").Value = 10598.5386
$ws.Range("L122").Value = 10101.9
$ws.Range("M122").Value = -8148.5386
$ws.Range("N122").Value = -15001.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 7372.1
$ws.Range("J75").Value = 10347
$ws.Range("L75").Value = 31041
$ws.Range("N75").Value = -33037
$ws.Range("H78").Value = 7372.1
$ws.Range("J78").Value = 10347
$ws.Range("L78").Value = 93123
$ws.Range("N78").Value = -103107
$ws.Range("H92").Value = 251.71428
$ws.Range("J92").Value = 270.75
$ws.Range("L92").Value = 812.25
$ws.Range("N92").Value = -3308.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 3667.1667
$ws.Range("J3").Value = 2800
$ws.Range("L3").Value = 2800
$ws.Range("N3").Value = -3032
$ws.Range("H70").Value = 2981
$ws.Range("I70").Value = 2198.7778
$ws.Range("J70").Value = 3986.7144
$ws.Range("K70").Value = 2198.7778
$ws.Range("L70").Value = 3986.7144
$ws.Range("M70").Value = -1928.7778
$ws.Range("N70").Value = -4526.7144
$ws.Range("H73").Value = 2981
$ws.Range("I73").Value = 2198.7778
$ws.Range("J73").Value = 3986.7144
$ws.Range("K73").Value = 2198.7778
$ws.Range("L73").Value = 3986.7144
$ws.Range("M73").Value = -1262.7778
$ws.Range("N73").Value = -5858.7144
$ws.Range("H113").Value = 3077
$ws.Range("I113").Value = 2967.1428
$ws.Range("J113").Value = 3333.3333
$ws.Range("K113").Value = 2967.1428
$ws.Range("L113").Value = 3333.3333
$ws.Range("M113").Value = -797.1428000000001
$ws.Range("N113").Value = -7673.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2071.5833
$ws.Range("I22").Value = 638.3333
$ws.Range("J22").Value = 2931.5334
$ws.Range("K22").Value = 638.3333
$ws.Range("L22").Value = 2931.5334
$ws.Range("M22").Value = -343.3333
$ws.Range("N22").Value = -3521.5334
$ws.Range("H27").Value = 2071.5833
$ws.Range("I27").Value = 638.3333
$ws.Range("J27").Value = 2931.5334
$ws.Range("K27").Value = 638.3333
$ws.Range("L27").Value = 2931.5334
$ws.Range("M27").Value = -531.3333
$ws.Range("N27").Value = -3145.5334
$ws.Range("H93").Value = 3746.5264
$ws.Range("I93").Value = 3827.9167
$ws.Range("J93").Value = 3607
$ws.Range("K93").Value = 3827.9167
$ws.Range("L93").Value = 3607
$ws.Range("M93").Value = -2579.9167
$ws.Range("N93").Value = -6103
$ws.Range("H122").Value = 3579.8076
$ws.Range("I122").Value = 2996.1667
$ws.Range("K122").Value = 8988.500100000001
$ws.Range("M122").Value = -6538.500100000001
$ws.Range("H132").Value = 2596.0862
$ws.Range("I132").Value = 2529.426
$ws.Range("K132").Value = 7588.278
$ws.Range("M132").Value = -5058.278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 73925.664
$ws.Range("J46").Value = 73925.664
$ws.Range("L46").Value = 73925.664
$ws.Range("N46").Value = -74387.664
$ws.Range("H96").Value = 4999.75
$ws.Range("J96").Value = 5000
$ws.Range("L96").Value = 5000
$ws.Range("N96").Value = -7746
$ws.Range("H100").Value = 1072.4884
$ws.Range("I100").Value = 1098.9395
$ws.Range("J100").Value = 985.2
$ws.Range("K100").Value = 2197.879
$ws.Range("L100").Value = 1970.4
$ws.Range("M100").Value = -1656.879
$ws.Range("N100").Value = -3052.4
$ws.Range("H107").Value = 727.3929000000001
$ws.Range("J107").Value = 854.8
$ws.Range("L107").Value = 2564.4
$ws.Range("N107").Value = -6404.4
$ws.Range("I126").Value = 2134
$ws.Range("J126").Value = 3069.7778
$ws.Range("K126").Value = 6402
$ws.Range("L126").Value = 9209.3334
$ws.Range("M126").Value = -3932
$ws.Range("N126").Value = -14149.3334
$ws.Range("H134").Value = 73925.664
$ws.Range("J134").Value = 73925.664
$ws.Range("L134").Value = 221776.992
$ws.Range("N134").Value = -226846.992
